# Update the NATMI LR-pair (Sema4d-Erbb2) results table with the recomputed TPM-based values.
# The table is expanded from a 3x4 (sending x target cluster) layout to a full 4x4 layout
# covering all four clusters (ECs, FAPs, MuSCs, Resolving-Mac) in both directions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4d"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.144900666666667
$ws.Range("H2").Value = 3.434702
$ws.Range("I2").Value = 0.02523133726002265
$ws.Range("J2").Value = 0.02523133726002265
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.020000333333333
$ws.Range("N2").Value = 9.060001
$ws.Range("O2").Value = 0.291481777372034
$ws.Range("P2").Value = 0.291481777372034
$ws.Range("Q2").Value = 3.457600394966889
$ws.Range("R2").Value = 31.118403554702
$ws.Range("S2").Value = 0.007354475030024627
$ws.Range("T2").Value = 0.007354475030024627

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4d"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.144900666666667
$ws.Range("H3").Value = 3.434702
$ws.Range("I3").Value = 0.02523133726002265
$ws.Range("J3").Value = 0.02523133726002265
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.076388666666666
$ws.Range("N3").Value = 12.229166
$ws.Range("O3").Value = 0.3934413518781783
$ws.Range("P3").Value = 0.3934413518781784
$ws.Range("Q3").Value = 4.667060102059111
$ws.Range("R3").Value = 42.003540918532
$ws.Range("S3").Value = 0.009927051441277561
$ws.Range("T3").Value = 0.009927051441277562

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4d"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.144900666666667
$ws.Range("H4").Value = 3.434702
$ws.Range("I4").Value = 0.02523133726002265
$ws.Range("J4").Value = 0.02523133726002265
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.229698
$ws.Range("N4").Value = 9.689094000000001
$ws.Range("O4").Value = 0.311721195201271
$ws.Range("P4").Value = 0.3117211952012711
$ws.Range("Q4").Value = 3.697683393332
$ws.Range("R4").Value = 33.279150539988
$ws.Range("S4").Value = 0.007865142607220621
$ws.Range("T4").Value = 0.007865142607220623

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema4d"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.144900666666667
$ws.Range("H5").Value = 3.434702
$ws.Range("I5").Value = 0.02523133726002265
$ws.Range("J5").Value = 0.02523133726002265
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03476766666666667
$ws.Range("N5").Value = 0.104303
$ws.Range("O5").Value = 0.003355675548516525
$ws.Range("P5").Value = 0.003355675548516525
$ws.Range("Q5").Value = 0.03980552474511111
$ws.Range("R5").Value = 0.358249722706
$ws.Range("S5").Value = 0.00008466818149983192
$ws.Range("T5").Value = 0.00008466818149983194

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4d"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.654659
$ws.Range("H6").Value = 4.963977
$ws.Range("I6").Value = 0.03646539869776051
$ws.Range("J6").Value = 0.03646539869776051
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.020000333333333
$ws.Range("N6").Value = 9.060001
$ws.Range("O6").Value = 0.291481777372034
$ws.Range("P6").Value = 0.291481777372034
$ws.Range("Q6").Value = 4.997070731552999
$ws.Range("R6").Value = 44.973636583977
$ws.Range("S6").Value = 0.01062899922500309
$ws.Range("T6").Value = 0.01062899922500309

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4d"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.654659
$ws.Range("H7").Value = 4.963977
$ws.Range("I7").Value = 0.03646539869776051
$ws.Range("J7").Value = 0.03646539869776051
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.076388666666666
$ws.Range("N7").Value = 12.229166
$ws.Range("O7").Value = 0.3934413518781783
$ws.Range("P7").Value = 0.3934413518781784
$ws.Range("Q7").Value = 6.745033194797999
$ws.Range("R7").Value = 60.70529875318199
$ws.Range("S7").Value = 0.01434699576042366
$ws.Range("T7").Value = 0.01434699576042366

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema4d"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.654659
$ws.Range("H8").Value = 4.963977
$ws.Range("I8").Value = 0.03646539869776051
$ws.Range("J8").Value = 0.03646539869776051
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.229698
$ws.Range("N8").Value = 9.689094000000001
$ws.Range("O8").Value = 0.311721195201271
$ws.Range("P8").Value = 0.3117211952012711
$ws.Range("Q8").Value = 5.344048862982
$ws.Range("R8").Value = 48.096439766838
$ws.Range("S8").Value = 0.01136703766555678
$ws.Range("T8").Value = 0.01136703766555678

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema4d"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.654659
$ws.Range("H9").Value = 4.963977
$ws.Range("I9").Value = 0.03646539869776051
$ws.Range("J9").Value = 0.03646539869776051
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.03476766666666667
$ws.Range("N9").Value = 0.104303
$ws.Range("O9").Value = 0.003355675548516525
$ws.Range("P9").Value = 0.003355675548516525
$ws.Range("Q9").Value = 0.057528632559
$ws.Range("R9").Value = 0.5177576930310001
$ws.Range("S9").Value = 0.0001223660467769813
$ws.Range("T9").Value = 0.0001223660467769813

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema4d"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.191626333333333
$ws.Range("H10").Value = 9.574878999999999
$ws.Range("I10").Value = 0.07033710676294723
$ws.Range("J10").Value = 0.07033710676294723
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.020000333333333
$ws.Range("N10").Value = 9.060001
$ws.Range("O10").Value = 0.291481777372034
$ws.Range("P10").Value = 0.291481777372034
$ws.Range("Q10").Value = 9.638712590542111
$ws.Range("R10").Value = 86.74841331487899
$ws.Range("S10").Value = 0.02050198489447037
$ws.Range("T10").Value = 0.02050198489447037

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema4d"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.191626333333333
$ws.Range("H11").Value = 9.574878999999999
$ws.Range("I11").Value = 0.07033710676294723
$ws.Range("J11").Value = 0.07033710676294723
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.076388666666666
$ws.Range("N11").Value = 12.229166
$ws.Range("O11").Value = 0.3934413518781783
$ws.Range("P11").Value = 0.3934413518781784
$ws.Range("Q11").Value = 13.01030941343489
$ws.Range("R11").Value = 117.092784720914
$ws.Range("S11").Value = 0.02767352637201372
$ws.Range("T11").Value = 0.02767352637201372

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema4d"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.191626333333333
$ws.Range("H12").Value = 9.574878999999999
$ws.Range("I12").Value = 0.07033710676294723
$ws.Range("J12").Value = 0.07033710676294723
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.229698
$ws.Range("N12").Value = 9.689094000000001
$ws.Range("O12").Value = 0.311721195201271
$ws.Range("P12").Value = 0.3117211952012711
$ws.Range("Q12").Value = 10.307989185514
$ws.Range("R12").Value = 92.771902669626
$ws.Range("S12").Value = 0.02192556698714531
$ws.Range("T12").Value = 0.02192556698714532

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema4d"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.191626333333333
$ws.Range("H13").Value = 9.574878999999999
$ws.Range("I13").Value = 0.07033710676294723
$ws.Range("J13").Value = 0.07033710676294723
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03476766666666667
$ws.Range("N13").Value = 0.104303
$ws.Range("O13").Value = 0.003355675548516525
$ws.Range("P13").Value = 0.003355675548516525
$ws.Range("Q13").Value = 0.1109654004818889
$ws.Range("R13").Value = 0.998688604337
$ws.Range("S13").Value = 0.0002360285093178183
$ws.Range("T13").Value = 0.0002360285093178184

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Sema4d"
$ws.Range("C14").Value = "Erbb2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 39.384953
$ws.Range("H14").Value = 118.154859
$ws.Range("I14").Value = 0.8679661572792696
$ws.Range("J14").Value = 0.8679661572792696
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.020000333333333
$ws.Range("N14").Value = 9.060001
$ws.Range("O14").Value = 0.291481777372034
$ws.Range("P14").Value = 0.291481777372034
$ws.Range("Q14").Value = 118.9425711883177
$ws.Range("R14").Value = 1070.483140694859
$ws.Range("S14").Value = 0.2529963182225359
$ws.Range("T14").Value = 0.2529963182225359

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Sema4d"
$ws.Range("C15").Value = "Erbb2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 39.384953
$ws.Range("H15").Value = 118.154859
$ws.Range("I15").Value = 0.8679661572792696
$ws.Range("J15").Value = 0.8679661572792696
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.076388666666666
$ws.Range("N15").Value = 12.229166
$ws.Range("O15").Value = 0.3934413518781783
$ws.Range("P15").Value = 0.3934413518781784
$ws.Range("Q15").Value = 160.5483760463993
$ws.Range("R15").Value = 1444.935384417594
$ws.Range("S15").Value = 0.3414937783044634
$ws.Range("T15").Value = 0.3414937783044634

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Sema4d"
$ws.Range("C16").Value = "Erbb2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 39.384953
$ws.Range("H16").Value = 118.154859
$ws.Range("I16").Value = 0.8679661572792696
$ws.Range("J16").Value = 0.8679661572792696
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.229698
$ws.Range("N16").Value = 9.689094000000001
$ws.Range("O16").Value = 0.311721195201271
$ws.Range("P16").Value = 0.3117211952012711
$ws.Range("Q16").Value = 127.201503934194
$ws.Range("R16").Value = 1144.813535407746
$ws.Range("S16").Value = 0.2705634479413483
$ws.Range("T16").Value = 0.2705634479413483

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema4d"
$ws.Range("C17").Value = "Erbb2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 39.384953
$ws.Range("H17").Value = 118.154859
$ws.Range("I17").Value = 0.8679661572792696
$ws.Range("J17").Value = 0.8679661572792696
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.03476766666666667
$ws.Range("N17").Value = 0.104303
$ws.Range("O17").Value = 0.003355675548516525
$ws.Range("P17").Value = 0.003355675548516525
$ws.Range("Q17").Value = 1.369322917586333
$ws.Range("R17").Value = 12.323906258277
$ws.Range("S17").Value = 0.002912612810921894
$ws.Range("T17").Value = 0.002912612810921894
